$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 167, shifting rows 167-173 down to 171-177.
$ws.Rows.Item(167).Insert()
$ws.Rows.Item(167).Insert()
$ws.Rows.Item(167).Insert()
$ws.Rows.Item(167).Insert()

# New row 167
$ws.Cells.Item(167,1).Value = 6
$ws.Cells.Item(167,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(167,3).Value = "Metropolitana"
$ws.Cells.Item(167,4).Value = 44509
$ws.Cells.Item(167,5).Value = 13
$ws.Cells.Item(167,6).Value = 100112026
$ws.Cells.Item(167,7).Value = "Haba"
$ws.Cells.Item(167,8).Value = "Sin especificar"
$ws.Cells.Item(167,9).Value = "Primera"
$ws.Cells.Item(167,10).Value = 320
$ws.Cells.Item(167,11).Value = 5000
$ws.Cells.Item(167,12).Value = 5000
$ws.Cells.Item(167,13).Value = 5000
$ws.Cells.Item(167,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(167,15).Value = "Región Metropolitana"
$ws.Cells.Item(167,16).Value = 200
$ws.Cells.Item(167,17).Value = 25
$ws.Cells.Item(167,18).Value = "Hortaliza"

# New row 168
$ws.Cells.Item(168,1).Value = 6
$ws.Cells.Item(168,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(168,3).Value = "Metropolitana"
$ws.Cells.Item(168,4).Value = 44509
$ws.Cells.Item(168,5).Value = 13
$ws.Cells.Item(168,6).Value = 100112026
$ws.Cells.Item(168,7).Value = "Haba"
$ws.Cells.Item(168,8).Value = "Sin especificar"
$ws.Cells.Item(168,9).Value = "Primera"
$ws.Cells.Item(168,10).Value = 1400
$ws.Cells.Item(168,11).Value = 5000
$ws.Cells.Item(168,12).Value = 6000
$ws.Cells.Item(168,13).Value = 5536
$ws.Cells.Item(168,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(168,15).Value = "Región del Maule"
$ws.Cells.Item(168,16).Value = 221
$ws.Cells.Item(168,17).Value = 25
$ws.Cells.Item(168,18).Value = "Hortaliza"

# New row 169
$ws.Cells.Item(169,1).Value = 6
$ws.Cells.Item(169,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(169,3).Value = "Metropolitana"
$ws.Cells.Item(169,4).Value = 44509
$ws.Cells.Item(169,5).Value = 13
$ws.Cells.Item(169,6).Value = 100112026
$ws.Cells.Item(169,7).Value = "Haba"
$ws.Cells.Item(169,8).Value = "Sin especificar"
$ws.Cells.Item(169,9).Value = "Segunda"
$ws.Cells.Item(169,10).Value = 460
$ws.Cells.Item(169,11).Value = 4000
$ws.Cells.Item(169,12).Value = 4000
$ws.Cells.Item(169,13).Value = 4000
$ws.Cells.Item(169,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(169,15).Value = "Región Metropolitana"
$ws.Cells.Item(169,16).Value = 160
$ws.Cells.Item(169,17).Value = 25
$ws.Cells.Item(169,18).Value = "Hortaliza"

# New row 170
$ws.Cells.Item(170,1).Value = 6
$ws.Cells.Item(170,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(170,3).Value = "Metropolitana"
$ws.Cells.Item(170,4).Value = 44509
$ws.Cells.Item(170,5).Value = 13
$ws.Cells.Item(170,6).Value = 100112026
$ws.Cells.Item(170,7).Value = "Haba"
$ws.Cells.Item(170,8).Value = "Sin especificar"
$ws.Cells.Item(170,9).Value = "Segunda"
$ws.Cells.Item(170,10).Value = 450
$ws.Cells.Item(170,11).Value = 4000
$ws.Cells.Item(170,12).Value = 4000
$ws.Cells.Item(170,13).Value = 4000
$ws.Cells.Item(170,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(170,15).Value = "Región del Maule"
$ws.Cells.Item(170,16).Value = 160
$ws.Cells.Item(170,17).Value = 25
$ws.Cells.Item(170,18).Value = "Hortaliza"
